# Generate Report for Handback
#
# The handback transform for e0b72c61-74a2-4c33-99eb-65cac4751436 failed
# because the returned file name did not match the expected handoff file
# name. Record the failure on the Overview sheet's Status column and add
# the detailed error message to the per-language "Error Detail" column,
# widening that column so the message is readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

# --- Overview sheet: update the Status for the e0b72c61... row ---------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = $statusText
$ovw.Range("F3").Value = $statusText

# --- zh-cn sheet: widen "Error Detail" column and record the error -----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $statusText
$zh.Columns.Item(16).ColumnWidth = 39.166666666666664
$zh.Range("P3").Value = "Handback file name: iqjxgyvv.hnq is different with handoff file name: e0b72c61-74a2-4c33-99eb-65cac4751436.4e00710dcccc00004b8cbbd2ebf5be0f93250884.zh-cn."

# --- de-de sheet: widen "Error Detail" column and record the error -----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $statusText
$de.Columns.Item(16).ColumnWidth = 39.166666666666664
$de.Range("P3").Value = "Handback file name: iqjxgyvv.hnq is different with handoff file name: e0b72c61-74a2-4c33-99eb-65cac4751436.4e00710dcccc00004b8cbbd2ebf5be0f93250884.de-de."
